# Feature: Adding new features in Province (#02)
# Appends a new school record ("LV") to the Province table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: nome | email | numero de salas | provincia
$ws.Range("A5").Value = "LV "
$ws.Range("B5").Value = "escolalv@gmail.com"
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = "Benguela"

# Turn the new email cell into a mailto hyperlink, matching the
# existing rows (B2:B4), and give it the same "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:escolalv@gmail.com")
$ws.Range("B5").Style = "Hyperlink"

# Move the active selection to the new last cell, like Excel does
# after data entry.
$ws.Range("D5").Select()
